# Add an extra "enhedstype" value on the "organisationenhed" sheet.
# Cell H3 currently re-uses the same enhedstype UUID as row 2 (the havn
# entry); give the Rådhuset row its own, distinct enhedstype UUID.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("organisationenhed")
$ws.Activate()

$ws.Range("H3").Value = "0034fa1f-b1ef-4764-8505-c5b9ca43aaa9"
[void]$ws.Range("H3").Select()
